# Adds a new column V ("07-10-2020") to the COVID19 active-cases sheet,
# mirroring the style of the existing date-header column U, and fills in
# the new daily figures for every state/UT row (2-36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header cell V1 : new date label, formatted like U1 -----------------
$headerCell = $ws.Cells.Item(1, 22)   # column V = 22

# Leading apostrophe forces the text "07-10-2020" to be stored as a plain
# string instead of being auto-parsed into a date serial number (this is
# exactly what Excel does for a user typed '07-10-2020).
$headerCell.Value = "'07-10-2020"

# Match the visual formatting used by the rest of the header row (N1:U1):
# bold font, centered/top aligned, thin box border all around.
$headerCell.HorizontalAlignment = -4108   # xlCenter
$headerCell.VerticalAlignment = -4160     # xlTop
$headerCell.Font.Bold = $true
$headerCell.Borders.LineStyle = 1         # xlContinuous (thin box border)

# ---- Data rows 2-36 : new numeric figures for 07-10-2020 ----------------
$newValues = @{
    2  = 180
    3  = 50776
    4  = 3022
    5  = 33047
    6  = 11420
    7  = 1492
    8  = 27238
    9  = 101
    10 = 22720
    11 = 4720
    12 = 16570
    13 = 11320
    14 = 3136
    15 = 13712
    16 = 10027
    17 = 115170
    18 = 87823
    19 = 1195
    20 = 18141
    21 = 247468
    22 = 2680
    23 = 2371
    24 = 261
    25 = 1185
    26 = 26846
    27 = 4522
    28 = 11982
    29 = 21294
    30 = 580
    31 = 45279
    32 = 26551
    33 = 4621
    34 = 8414
    35 = 44031
    36 = 27988
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 22).Value = $newValues[$row]
}
